$d = $word.ActiveDocument

# Insert a new, empty paragraph right after the last paragraph in the body
# ("Complete all") and before the sectPr. The new paragraph keeps the
# "ListParagraph" style and bold character formatting (b/bCs) on its
# paragraph mark, but carries no numbering and no text run.
$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)

$newParagraphXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>'

[void]$insertionPoint.InsertXML($newParagraphXml)
